$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.841.54"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.02"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.55"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.22"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.15"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.79"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.34"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.849.73"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.33"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.79"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.284.45"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.803"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.782.25"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("E44").Value = "  -6.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.60"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.01"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.57"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -0.17%  "
